$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 3's formatting (styles only) down into row 4 so the new row
# inherits the same look (e.g. text-formatted date/number columns),
# then fill in the new test-case row's values explicitly.
$ws.Range("A3:L3").Copy()
$ws.Range("A4:L4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A4").Value = "checkFilters"
$ws.Range("B4").Value = "Bangalore"
$ws.Range("C4").Value = "Tokyo"
$ws.Range("D4").Value = "Bangalore, IN - Kempegowda International Airport (BLR)"
$ws.Range("E4").Value = "Tokyo, JP - Narita (NRT)"
$ws.Range("F4").Value = "2/12/2017"
$ws.Range("G4").Value = "Mr Akash Sharma"
$ws.Range("H4").Value = "1/Jan/1993"
$ws.Range("I4").Value = 987654321
$ws.Range("J4").Value = "Business"
$ws.Range("K4").Value = 987654321
$ws.Range("L4").Value = "testUser2@gmail.com"

# Match the new selection state recorded in the sheet view.
$ws.Range("M4").Select()
